$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the worker rows: move WILBERTO / ALBERTO up to rows 16-17,
# and push RONALD / HECTOR down to rows 18-19. Also update HECTOR's
# "Valor Mora" (G19) to the new value from the updated EC database.

# Row 16 -> WILBERTO ANIBAL TEHERAN GUTIERREZ
$ws.Range("C16").Value = "1050946110"
$ws.Range("D16").Value = "WILBERTO ANIBAL TEHERAN GUTIERREZ"
$ws.Range("E16").Value = "1808"
$ws.Range("F16").Value = 1240
$ws.Range("G16").Value = 930000

# Row 17 -> ALBERTO LUIS CANABAL MARRUGO
$ws.Range("C17").Value = "1051443479"
$ws.Range("D17").Value = "ALBERTO LUIS CANABAL MARRUGO"
$ws.Range("E17").Value = "1808"
$ws.Range("F17").Value = 1240
$ws.Range("G17").Value = 930000

# Row 18 -> RONALD DE JESUS PUELLO BARRIOS
$ws.Range("C18").Value = "1050962190"
$ws.Range("D18").Value = "RONALD DE JESUS PUELLO BARRIOS"
$ws.Range("E18").Value = "2201"
$ws.Range("F18").Value = 16867
$ws.Range("G18").Value = 1150000

# Row 19 -> HECTOR DAVID SALGADO ARTEAGA (updated Valor Mora)
$ws.Range("C19").Value = "1050954380"
$ws.Range("D19").Value = "HECTOR DAVID SALGADO ARTEAGA"
$ws.Range("E19").Value = "2205"
$ws.Range("F19").Value = 38333
$ws.Range("G19").Value = 1724688
